$d = $word.ActiveDocument

function Insert-ParagraphAfterText {
    param([string]$AnchorText, [string]$NewText)
    $rng = $d.Content
    $found = $rng.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $AnchorText"
    }
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $countRng = $d.Range(0, $rng.Start)
    $idx = $countRng.Paragraphs.Count
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $NewText
    return $idx + 1
}

# 1. After "List of sports" paragraph, add a new ListParagraph (ilvl 1) bullet.
Insert-ParagraphAfterText "List of sports" "Category:National Basketball Association teams" | Out-Null

# 2. After "Geography" paragraph, add a new ListParagraph (ilvl 1) bullet "Continents".
Insert-ParagraphAfterText "Geography" "Continents" | Out-Null

# 3. After "Category:Logic puzzles" paragraph, insert a new paragraph "List of television
#    programs by name" (ListParagraph, ilvl 1) and move the _GoBack bookmark onto it.
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Category:Logic puzzles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Anchor text not found: Category:Logic puzzles" }
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Category:Logic puzzles</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>List of television programs by name</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rng3.InsertXML($xml3)

# 4. After the "Look into ... (Not working yet!):" paragraph, add:
#      - an empty paragraph
#      - "List of 20th Century Fox films" (with "th" superscripted)
#      - "List of television programs by name" with a tab stop at 2430 twips, split into
#        3 runs ("List of television program" / "s" / " by name")
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("(Not working yet!):", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "Anchor text not found: (Not working yet!):" }
$rng4.Collapse(0)
$rng4.InsertParagraphAfter()
$countRng4 = $d.Range(0, $rng4.Start)
$idx4 = $countRng4.Paragraphs.Count
# $idx4 now points at the just-created empty paragraph; leave it blank (<w:p/>).

# Insert the "List of 20th Century Fox films" paragraph right after the blank one.
$blankPara = $d.Paragraphs.Item($idx4 + 1)
$blankPara.Range.InsertParagraphAfter()
$filmsPara = $d.Paragraphs.Item($idx4 + 2)
$filmsPara.Range.Text = "List of 20th Century Fox films"
$pStart = $filmsPara.Range.Start
$thStart = $pStart + 10
$thEnd = $thStart + 2
$thRange = $d.Range($thStart, $thEnd)
$thRange.Font.Superscript = $true

# Insert the final "List of television programs by name" paragraph (with tab stop and
# 3 separate runs) right after the films paragraph.
$filmsPara2 = $d.Paragraphs.Item($idx4 + 2)
$filmsPara2.Range.InsertParagraphAfter()
$tvPara = $d.Paragraphs.Item($idx4 + 3)
$tvXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="2430"/></w:tabs></w:pPr><w:r><w:t>List of television program</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> by name</w:t></w:r></w:p>'
$tvPara.Range.InsertXML($tvXml)
